$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column A labels (rows 2-15) with new shared-string text
$ws.Range("A2").Value = "Global climate scheme (GCS)"
$ws.Range("A3").Value = "National tax on millionaires funding public services"
$ws.Range("A4").Value = "Global tax on millionaires funding low-income countries"
$ws.Range("A5").Value = "Preferred share of global wealth`ntax for low-income countries: ≥ 30%"
$ws.Range("A6").Value = "[Country]'s foreign aid should be increased*"
$ws.Range("A7").Value = "High-income countries contributing $100 billion per year`nto help low-income countries adapt to climate change"
$ws.Range("A8").Value = "High-income countries funding renewable`nenergy in low-income countries"
$ws.Range("A9").Value = "Payments from high-income countries to compensate`nlow-income countries for climate damages"
$ws.Range("A10").Value = "Cancellation of low-income countries' public debt"
$ws.Range("A11").Value = "Democratise international institutions (UN, IMF) by making`na country's voting right proportional to its population"
$ws.Range("A12").Value = "Removing tariffs on imports from low-income countries"
$ws.Range("A13").Value = "A minimum wage in all countries`nat 50% of local median wage"
$ws.Range("A14").Value = "Fight tax evasion by creating a global financial`nregister to record ownership of all assets"
$ws.Range("A15").Value = "A maximum wealth limit of $10 billion`n(US) / €100 million (Eu) for each human"

# Update numeric data B:G for rows 2-15
$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.542040004729187
$arr[0,1] = 0.757320866764204
$arr[0,2] = 0.802845995450502
$arr[0,3] = 0.712681465751731
$arr[0,4] = 0.809917713113721
$arr[0,5] = 0.74106127773703
$ws.Range("B2:G2").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.619576204293734
$arr[0,1] = 0.758375208534753
$arr[0,2] = 0.697490019207693
$arr[0,3] = 0.785605734377313
$arr[0,4] = 0.788109233298969
$arr[0,5] = 0.774030617727549
$ws.Range("B3:G3").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.580974048697037
$arr[0,1] = 0.714793319488843
$arr[0,2] = 0.694319177328758
$arr[0,3] = 0.723244532725875
$arr[0,4] = 0.776824002261202
$arr[0,5] = 0.708989998773392
$ws.Range("B4:G4").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.502686637940573
$arr[0,1] = 0.535405079823794
$arr[0,2] = 0.53495530999455
$arr[0,3] = 0.504405401186565
$arr[0,4] = 0.57374820373221
$arr[0,5] = 0.54094700051697
$ws.Range("B5:G5").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.598187005678365
$arr[0,1] = 0.636739028674537
$arr[0,2] = 0.627588384954585
$arr[0,3] = 0.676807601959896
$arr[0,4] = 0.691059073949057
$arr[0,5] = 0.560844998477962
$ws.Range("B6:G6").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.45011340420921
$arr[0,1] = 0.580303842050924
$arr[0,2] = 0.554844714670846
$arr[0,3] = 0.597860385545241
$arr[0,4] = 0.624386789670735
$arr[0,5] = 0.54091224731223
$ws.Range("B7:G7").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.53081419263325
$arr[0,1] = 0.646086050975591
$arr[0,2] = 0.617631776703308
$arr[0,3] = 0.661981009518735
$arr[0,4] = 0.679241918379181
$arr[0,5] = 0.624160044931323
$ws.Range("B8:G8").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.409115346811486
$arr[0,1] = 0.535873014459192
$arr[0,2] = 0.522108726136302
$arr[0,3] = 0.528761543405907
$arr[0,4] = 0.619349694288149
$arr[0,5] = 0.510164860189199
$ws.Range("B9:G9").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.308391076954102
$arr[0,1] = 0.365731640324952
$arr[0,2] = 0.362994387715357
$arr[0,3] = 0.298641910128985
$arr[0,4] = 0.450925576392201
$arr[0,5] = 0.404799933538172
$ws.Range("B10:G10").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.3441053797594
$arr[0,1] = 0.443594773467438
$arr[0,2] = 0.436219325369461
$arr[0,3] = 0.428402025624624
$arr[0,4] = 0.517519698549576
$arr[0,5] = 0.432189204736391
$ws.Range("B11:G11").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.390484602228333
$arr[0,1] = 0.49094991581821
$arr[0,2] = 0.386302468418039
$arr[0,3] = 0.512592712201629
$arr[0,4] = 0.502893401333253
$arr[0,5] = 0.540728143603005
$ws.Range("B12:G12").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.41536013985638
$arr[0,1] = 0.547947048183082
$arr[0,2] = 0.541679032092289
$arr[0,3] = 0.542087544034292
$arr[0,4] = 0.60768384976524
$arr[0,5] = 0.531644320070783
$ws.Range("B13:G13").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.441485472151688
$arr[0,1] = 0.699336514515305
$arr[0,2] = 0.730152208937297
$arr[0,3] = 0.702462385725501
$arr[0,4] = 0.718628182752579
$arr[0,5] = 0.6490870831924
$ws.Range("B14:G14").Value = $arr

$arr = New-Object 'object[,]' 1,6
$arr[0,0] = 0.335948876827959
$arr[0,1] = 0.45314379468064
$arr[0,2] = 0.395498925733041
$arr[0,3] = 0.479942664267705
$arr[0,4] = 0.441371326747748
$arr[0,5] = 0.496306674811346
$ws.Range("B15:G15").Value = $arr
